# Generate Report for Handoff
# Adds a new file entry (a4a7b6de-7656-4320-ae58-6c9ffbac640a.md) as row 3
# on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newMd = "a4a7b6de-7656-4320-ae58-6c9ffbac640a.md"
$newMdPath = "e2e\a4a7b6de-7656-4320-ae58-6c9ffbac640a.md"
$newMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0974f356a11ef705bf04685e24cb29ec6603873e/e2e/a4a7b6de-7656-4320-ae58-6c9ffbac640a.md"

$hoDateTime = "2016-09-02 14:49:35"
$zhTargetFile = "a4a7b6de-7656-4320-ae58-6c9ffbac640a.5fad6ea8ff1e89253400eb90302be13ff86c1176.zh-cn.xlf"
$zhHoDatetime = "2016-09-02 14:49:30"
$deTargetFile = "a4a7b6de-7656-4320-ae58-6c9ffbac640a.5fad6ea8ff1e89253400eb90302be13ff86c1176.de-de.xlf"
$deHoDatetime = $hoDateTime

$linkColor = 15570276  # cornflower blue (BGR order) matching existing hyperlink style

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", $newMdPath)
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = $linkColor
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $hoDateTime
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newMdUrl, "", "", $newMd)
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = $linkColor
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhTargetFile
$wsZhCn.Range("H3").Value = $zhHoDatetime
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newMdUrl, "", "", $newMd)
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = $linkColor
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $deTargetFile
$wsDeDe.Range("H3").Value = $deHoDatetime
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

Write-Host "Report rows appended for handoff."
